# Generate Report for Archive
# Update the localization status of both e2e files from "Ready for handoff"
# to "In Translation" on every sheet that surfaces it (Overview's per-locale
# status columns, and the Status column on each locale's detail sheet), then
# resize the affected status columns to their new best-fit width.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn (E) and de-de (F) status columns, rows 2-3
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus

# zh-cn detail sheet: Status column (C), rows 2-3
$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus

# de-de detail sheet: Status column (C), rows 2-3
$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# The status text got shorter, so the columns that display it shrink to the
# new best-fit width.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
